$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92; existing rows 92..170 shift down to 93..171.
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new record.
$ws.Cells.Item(92, 1).Value = 5
$ws.Cells.Item(92, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(92, 3).Value = "Maule"
$ws.Cells.Item(92, 4).Value = 44790
$ws.Cells.Item(92, 5).Value = 7
$ws.Cells.Item(92, 6).Value = 100112031
$ws.Cells.Item(92, 7).Value = "Poroto verde"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 100
$ws.Cells.Item(92, 11).Value = 35000
$ws.Cells.Item(92, 12).Value = 35000
$ws.Cells.Item(92, 13).Value = 35000
$ws.Cells.Item(92, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(92, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(92, 16).Value = 1400
$ws.Cells.Item(92, 17).Value = 25
$ws.Cells.Item(92, 18).Value = "Hortaliza"
